$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.115.10'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').Value = '3.733.93'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.06'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.00'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').Value = '3.731.47'
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.166'
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.39'
$ws.Range('E11').Value = '  +3.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.06'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '4.355.63'
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('D16').Value = '3.735.26'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('D17').Value = '69.156.88'
$ws.Range('E17').Value = '  +2.51%  '
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('E19').Value = '  -0.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.02'
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.96'
$ws.Range('E21').Value = '  +18.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '493.52'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.725'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('E24').Value = '  +9.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.70'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.31'
$ws.Range('E26').Value = '  +0.85%  '
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.14'
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.00'
$ws.Range('E30').Value = '  +2.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.50'
$ws.Range('E31').Value = '  +6.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.05'
$ws.Range('E32').Value = '  +4.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.58'
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('D34').Value = '3.879.27'
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.108'
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').Value = '3.666.97'
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.86'
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('E40').Value = '  +2.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.323'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.99'
$ws.Range('E42').Value = '  +6.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '431.47'
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.64'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.46'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.20'
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.55'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').Value = '2.778.63'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('E51').Value = '  +0.62%  '
